$d = $word.ActiveDocument
$d.Content.Find.Execute("url = http://localhost:3000/", $true, $false, $false, $false, $false, $true, 1, $false, "url = http://localhost:8080/", 2)
